$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1269.6666
$ws.Range("I17").Value = 812
$ws.Range("J17").Value = 1498.5
$ws.Range("K17").Value = 2436
$ws.Range("L17").Value = 4495.5
$ws.Range("M17").Value = -2268
$ws.Range("N17").Value = -4831.5
$ws.Range("H132").Value = 9829.333000000001
$ws.Range("I132").Value = 9359.208000000001
$ws.Range("K132").Value = 28077.624
$ws.Range("M132").Value = -25547.624
$ws.Range("H138").Value = 3267.65
$ws.Range("I138").Value = 897.5
$ws.Range("K138").Value = 2692.5
$ws.Range("M138").Value = 2447.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1001.7273
$ws.Range("I2").Value = 902.1111
$ws.Range("K2").Value = 902.1111
$ws.Range("M2").Value = -789.1111
$ws.Range("H32").Value = 4344.8184
$ws.Range("I32").Value = 4313.619
$ws.Range("K32").Value = 4313.619
$ws.Range("M32").Value = -4026.619
$ws.Range("H74").Value = 2549.077
$ws.Range("I74").Value = 2549.077
$ws.Range("K74").Value = 2549.077
$ws.Range("M74").Value = -1675.077
$ws.Range("H77").Value = 2549.077
$ws.Range("I77").Value = 2549.077
$ws.Range("K77").Value = 12745.385
$ws.Range("M77").Value = -8377.385000000002
$ws.Range("H116").Value = 1001.7273
$ws.Range("I116").Value = 902.1111
$ws.Range("K116").Value = 902.1111
$ws.Range("M116").Value = 1391.8889
$ws.Range("H132").Value = 5406.75
$ws.Range("I132").Value = 5695.5
$ws.Range("K132").Value = 17086.5
$ws.Range("M132").Value = -14556.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1001.7273
$ws.Range("I3").Value = 902.1111
$ws.Range("K3").Value = 902.1111
$ws.Range("M3").Value = -788.1111
$ws.Range("H22").Value = 1319
$ws.Range("I22").Value = 1347.8182
$ws.Range("J22").Value = 1002
$ws.Range("K22").Value = 1347.8182
$ws.Range("L22").Value = 1002
$ws.Range("M22").Value = -1174.8182
$ws.Range("N22").Value = -1348
$ws.Range("H29").Value = 627.7143
$ws.Range("J29").Value = 1050
$ws.Range("L29").Value = 1050
$ws.Range("N29").Value = -1628
$ws.Range("H36").Value = 76.333336
$ws.Range("I36").Value = 64.5
$ws.Range("J36").Value = 100
$ws.Range("K36").Value = 64.5
$ws.Range("L36").Value = 100
$ws.Range("M36").Value = 469.5
$ws.Range("N36").Value = -1168
$ws.Range("H105").Value = 2045.7142
$ws.Range("I105").Value = 1963.8
$ws.Range("J105").Value = 2250.5
$ws.Range("K105").Value = 1963.8
$ws.Range("L105").Value = 2250.5
$ws.Range("M105").Value = -216.8
$ws.Range("N105").Value = -5744.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5962.48
$ws.Range("I31").Value = 2641.0833
$ws.Range("J31").Value = 9028.385
$ws.Range("K31").Value = 2641.0833
$ws.Range("L31").Value = 9028.385
$ws.Range("M31").Value = -2346.0833
$ws.Range("N31").Value = -9618.385
$ws.Range("H34").Value = 5962.48
$ws.Range("I34").Value = 2641.0833
$ws.Range("J34").Value = 9028.385
$ws.Range("K34").Value = 2641.0833
$ws.Range("L34").Value = 9028.385
$ws.Range("M34").Value = -2439.0833
$ws.Range("N34").Value = -9432.385
$ws.Range("H88").Value = 9999.333000000001
$ws.Range("J88").Value = 9999.333000000001
$ws.Range("L88").Value = 9999.333000000001
$ws.Range("N88").Value = -10811.333
$ws.Range("H91").Value = 9999.333000000001
$ws.Range("J91").Value = 9999.333000000001
$ws.Range("L91").Value = 9999.333000000001
$ws.Range("N91").Value = -12807.333
$ws.Range("H96").Value = 10072
$ws.Range("J96").Value = 10072
$ws.Range("L96").Value = 10072
$ws.Range("N96").Value = -15564
$ws.Range("H105").Value = 2761.3076
$ws.Range("I105").Value = 2649.6667
$ws.Range("J105").Value = 2857
$ws.Range("K105").Value = 2649.6667
$ws.Range("L105").Value = 2857
$ws.Range("M105").Value = -902.6667000000002
$ws.Range("N105").Value = -6351
$ws.Range("H132").Value = 994.6
$ws.Range("I132").Value = 994.6
$ws.Range("K132").Value = 2983.8
$ws.Range("M132").Value = -453.8000000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1279.6666
$ws.Range("I18").Value = 1279.6666
$ws.Range("K18").Value = 3838.9998
$ws.Range("M18").Value = -3669.9998
$ws.Range("H63").Value = 1999
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 1999
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 5997
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -7495
$ws.Range("H66").Value = 1999
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 1999
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 17991
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -25479
$ws.Range("H92").Value = 721.75
$ws.Range("J92").Value = 900
$ws.Range("L92").Value = 2700
$ws.Range("N92").Value = -5196
$ws.Range("H107").Value = 504.16666
$ws.Range("I107").Value = 366.66666
$ws.Range("K107").Value = 1099.99998
$ws.Range("M107").Value = 820.0000199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H136").Value = 3813.6
$ws.Range("J136").Value = 4987.5
$ws.Range("L136").Value = 14962.5
$ws.Range("N136").Value = -20062.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9760
$ws.Range("H122").Value = 3217.818
$ws.Range("I122").Value = 1771
$ws.Range("K122").Value = 5313
$ws.Range("M122").Value = -2863
$ws.Range("H126").Value = 6494.2
$ws.Range("I126").Value = 4581.2
$ws.Range("J126").Value = 7450.7
$ws.Range("K126").Value = 13743.6
$ws.Range("L126").Value = 22352.1
$ws.Range("M126").Value = -11273.6
$ws.Range("N126").Value = -27292.1
$ws.Range("H136").Value = 3493.24
$ws.Range("I136").Value = 2352.2856
$ws.Range("J136").Value = 4945.364
$ws.Range("K136").Value = 7056.8568
$ws.Range("L136").Value = 14836.092
$ws.Range("M136").Value = -4506.8568
$ws.Range("N136").Value = -19936.092
